$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Leave")
$ws.Name = "Trip"
